# Apply the updated crypto price/volume values from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.557.01'
$ws.Range('E2').Value = '  +3.97%  '
$ws.Range('D3').Value = '1.814.39'
$ws.Range('E3').Value = '  +5.32%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '334.51'
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3818'
$ws.Range('E7').Value = '  +2.21%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3512'
$ws.Range('E8').Value = '  +4.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '49.46'
$ws.Range('E9').Value = '  +2.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.226'
$ws.Range('E10').Value = '  +3.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07677'
$ws.Range('E11').Value = '  +4.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.004'
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.14'
$ws.Range('E13').Value = '  +9.72%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.585'
$ws.Range('E14').Value = '  +3.36%  '
$ws.Range('D15').Value = '1.820.12'
$ws.Range('E15').Value = '  +5.64%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.136'
$ws.Range('E16').Value = '  +1.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001116'
$ws.Range('E17').Value = '  +4.29%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06700'
$ws.Range('E18').Value = '  +0.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '86.65'
$ws.Range('E19').Value = '  +5.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.48'
$ws.Range('E21').Value = '  +5.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.505'
$ws.Range('E22').Value = '  +6.61%  '
$ws.Range('B23').Value = 'WrappedBTC'
$ws.Range('C23').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D23').Value = '27.561.32'
$ws.Range('E23').Value = '  +4.31%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.03'
$ws.Range('E24').Value = '  +1.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.458'
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.636'
$ws.Range('E26').Value = '  +10.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.97'
$ws.Range('E27').Value = '  +13.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.471'
$ws.Range('E28').Value = '  +5.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '151.34'
$ws.Range('E29').Value = '  -1.14%  '
$ws.Range('D30').Value = '2.026.92'
$ws.Range('E30').Value = '  +5.99%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '134.65'
$ws.Range('E31').Value = '  +3.16%  '
$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.099'
$ws.Range('E32').Value = '  -1.06%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.241'
$ws.Range('E33').Value = '  +4.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '13.71'
$ws.Range('E34').Value = '  +8.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.08717'
$ws.Range('E35').Value = '  +1.26%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.690'
$ws.Range('E36').Value = '  -0.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.561'
$ws.Range('E37').Value = '  +3.56%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6937'
$ws.Range('E38').Value = '  +12.49%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '9.034'
$ws.Range('E39').Value = '  +6.60%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2244'
$ws.Range('E40').Value = '  +4.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.06460'
$ws.Range('E41').Value = '  +4.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.02380'
$ws.Range('E42').Value = '  +2.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.297'
$ws.Range('E43').Value = '  +5.96%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.76'
$ws.Range('E44').Value = '  +5.58%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6517'
$ws.Range('E45').Value = '  +9.39%  '
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.868'
$ws.Range('E47').Value = '  -0.82%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.164'
$ws.Range('E48').Value = '  +6.24%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '132.50'
$ws.Range('E49').Value = '  +3.62%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07292'
$ws.Range('E50').Value = '  +1.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '80.76'
$ws.Range('E51').Value = '  +5.38%  '
